# Apply cryptos.xlsx price/volume/coin updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.958.26"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.212.74"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D5").Value = "'241.82"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").Value = "'73.36"
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "'43.42"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("D11").Value = "'0.0954"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "'7.08"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "2.551.71"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").Value = "'14.21"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "'0.842"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "2.226.68"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "41.822.67"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("E19").Value = "  +9.94%  "
$ws.Range("D20").Value = "'72.43"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "'6.14"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'10.45"
$ws.Range("E22").Value = "  +17.30%  "
$ws.Range("D23").Value = "'229.03"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "'2.08"
$ws.Range("E24").Value = "  -6.90%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'11.48"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("D27").Value = "'3.58"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "'2.15"
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("D30").Value = "'166.81"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").Value = "'20.53"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'5.56"
$ws.Range("E32").Value = "  +6.11%  "
$ws.Range("D33").Value = "'0.0792"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'28.81"
$ws.Range("E35").Value = "  -4.74%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  -6.40%  "
$ws.Range("D37").Value = "'4.26"
$ws.Range("E37").Value = "  -6.05%  "
$ws.Range("D38").Value = "'0.0301"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "'12.80"
$ws.Range("E39").Value = "  -4.87%  "
$ws.Range("D40").Value = "'65.93"
$ws.Range("E40").Value = "  +6.47%  "
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("D42").Value = "'5.61"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").Value = "'0.200"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").Value = "'8.71"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'103.66"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("D46").Value = "'0.100"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").Value = "'2.41"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("D48").Value = "'1.11"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").Value = "'2.71"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "2.422.67"
$ws.Range("E51").Value = "  -1.35%  "
